$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Sema3a"
$ws.Cells.Item(2, 3).Value = "Plxna1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.397441
$ws.Cells.Item(2, 8).Value = 4.192323
$ws.Cells.Item(2, 9).Value = 0.6676161521996591
$ws.Cells.Item(2, 10).Value = 0.6676161521996592
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 7.499519333333335
$ws.Cells.Item(2, 14).Value = 22.498558
$ws.Cells.Item(2, 15).Value = 0.2063210495448665
$ws.Cells.Item(2, 16).Value = 0.2063210495448665
$ws.Cells.Item(2, 17).Value = 10.48013579669267
$ws.Cells.Item(2, 18).Value = 94.32122217023401
$ws.Cells.Item(2, 19).Value = 0.137743265214939
$ws.Cells.Item(2, 20).Value = 0.137743265214939

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Sema3a"
$ws.Cells.Item(3, 3).Value = "Plxna1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.397441
$ws.Cells.Item(3, 8).Value = 4.192323
$ws.Cells.Item(3, 9).Value = 0.6676161521996591
$ws.Cells.Item(3, 10).Value = 0.6676161521996592
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 14.856814
$ws.Cells.Item(3, 14).Value = 44.570442
$ws.Cells.Item(3, 15).Value = 0.4087293226578609
$ws.Cells.Item(3, 16).Value = 0.4087293226578609
$ws.Cells.Item(3, 17).Value = 20.761521012974
$ws.Cells.Item(3, 18).Value = 186.853689116766
$ws.Cells.Item(3, 19).Value = 0.272874297684014
$ws.Cells.Item(3, 20).Value = 0.2728742976840141

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Sema3a"
$ws.Cells.Item(4, 3).Value = "Plxna1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.397441
$ws.Cells.Item(4, 8).Value = 4.192323
$ws.Cells.Item(4, 9).Value = 0.6676161521996591
$ws.Cells.Item(4, 10).Value = 0.6676161521996592
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 13.992451
$ws.Cells.Item(4, 14).Value = 41.977353
$ws.Cells.Item(4, 15).Value = 0.3849496277972726
$ws.Cells.Item(4, 16).Value = 0.3849496277972725
$ws.Cells.Item(4, 17).Value = 19.553624717891
$ws.Cells.Item(4, 18).Value = 175.982622461019
$ws.Cells.Item(4, 19).Value = 0.2569985893007061
$ws.Cells.Item(4, 20).Value = 0.2569985893007061

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Sema3a"
$ws.Cells.Item(5, 3).Value = "Plxna1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.2347746666666667
$ws.Cells.Item(5, 8).Value = 0.7043240000000001
$ws.Cells.Item(5, 9).Value = 0.112161700990566
$ws.Cells.Item(5, 10).Value = 0.112161700990566
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 7.499519333333335
$ws.Cells.Item(5, 14).Value = 22.498558
$ws.Cells.Item(5, 15).Value = 0.2063210495448665
$ws.Cells.Item(5, 16).Value = 0.2063210495448665
$ws.Cells.Item(5, 17).Value = 1.760697151643556
$ws.Cells.Item(5, 18).Value = 15.846274364792
$ws.Cells.Item(5, 19).Value = 0.02314131986711108
$ws.Cells.Item(5, 20).Value = 0.02314131986711107

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Sema3a"
$ws.Cells.Item(6, 3).Value = "Plxna1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.2347746666666667
$ws.Cells.Item(6, 8).Value = 0.7043240000000001
$ws.Cells.Item(6, 9).Value = 0.112161700990566
$ws.Cells.Item(6, 10).Value = 0.112161700990566
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 14.856814
$ws.Cells.Item(6, 14).Value = 44.570442
$ws.Cells.Item(6, 15).Value = 0.4087293226578609
$ws.Cells.Item(6, 16).Value = 0.4087293226578609
$ws.Cells.Item(6, 17).Value = 3.488003554578667
$ws.Cells.Item(6, 18).Value = 31.392031991208
$ws.Cells.Item(6, 19).Value = 0.04584377607402757
$ws.Cells.Item(6, 20).Value = 0.04584377607402758

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Sema3a"
$ws.Cells.Item(7, 3).Value = "Plxna1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.2347746666666667
$ws.Cells.Item(7, 8).Value = 0.7043240000000001
$ws.Cells.Item(7, 9).Value = 0.112161700990566
$ws.Cells.Item(7, 10).Value = 0.112161700990566
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 13.992451
$ws.Cells.Item(7, 14).Value = 41.977353
$ws.Cells.Item(7, 15).Value = 0.3849496277972726
$ws.Cells.Item(7, 16).Value = 0.3849496277972725
$ws.Cells.Item(7, 17).Value = 3.285073019374667
$ws.Cells.Item(7, 18).Value = 29.565657174372
$ws.Cells.Item(7, 19).Value = 0.04317660504942737
$ws.Cells.Item(7, 20).Value = 0.04317660504942737

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Sema3a"
$ws.Cells.Item(8, 3).Value = "Plxna1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.4609646666666667
$ws.Cells.Item(8, 8).Value = 1.382894
$ws.Cells.Item(8, 9).Value = 0.2202221468097748
$ws.Cells.Item(8, 10).Value = 0.2202221468097748
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 7.499519333333335
$ws.Cells.Item(8, 14).Value = 22.498558
$ws.Cells.Item(8, 15).Value = 0.2063210495448665
$ws.Cells.Item(8, 16).Value = 0.2063210495448665
$ws.Cells.Item(8, 17).Value = 3.457013429650223
$ws.Cells.Item(8, 18).Value = 31.113120866852
$ws.Cells.Item(8, 19).Value = 0.04543646446281641
$ws.Cells.Item(8, 20).Value = 0.04543646446281641

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Sema3a"
$ws.Cells.Item(9, 3).Value = "Plxna1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.4609646666666667
$ws.Cells.Item(9, 8).Value = 1.382894
$ws.Cells.Item(9, 9).Value = 0.2202221468097748
$ws.Cells.Item(9, 10).Value = 0.2202221468097748
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 14.856814
$ws.Cells.Item(9, 14).Value = 44.570442
$ws.Cells.Item(9, 15).Value = 0.4087293226578609
$ws.Cells.Item(9, 16).Value = 0.4087293226578609
$ws.Cells.Item(9, 17).Value = 6.848466313238667
$ws.Cells.Item(9, 18).Value = 61.636196819148
$ws.Cells.Item(9, 19).Value = 0.09001124889981925
$ws.Cells.Item(9, 20).Value = 0.09001124889981926

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Sema3a"
$ws.Cells.Item(10, 3).Value = "Plxna1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.4609646666666667
$ws.Cells.Item(10, 8).Value = 1.382894
$ws.Cells.Item(10, 9).Value = 0.2202221468097748
$ws.Cells.Item(10, 10).Value = 0.2202221468097748
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 13.992451
$ws.Cells.Item(10, 14).Value = 41.977353
$ws.Cells.Item(10, 15).Value = 0.3849496277972726
$ws.Cells.Item(10, 16).Value = 0.3849496277972725
$ws.Cells.Item(10, 17).Value = 6.450025511064667
$ws.Cells.Item(10, 18).Value = 58.050229599582
$ws.Cells.Item(10, 19).Value = 0.08477443344713913
$ws.Cells.Item(10, 20).Value = 0.08477443344713913
